$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update labels in column A
$ws.Range("A2").Value = "Wafer 1"
$ws.Range("A3").Value = "Wafer 2"

# Update mean values in column B
$ws.Range("B2").Value = 0.4832801460336079
$ws.Range("B3").Value = 0.4822122652746012

# Remove row 4 entirely (previously "Blue" / 28.16764333333333)
$ws.Rows.Item(4).Delete()
